$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row for "shexview" right after row 94 (becomes the new row 95) ---
$ws.Rows.Item(95).Insert() | Out-Null

# Copy cell formatting from the row above (row 94, "setpoint") so the new row
# matches the existing look & feel of the table instead of getting bare defaults.
for ($c = 1; $c -le 8; $c++) {
    $src = $ws.Cells.Item(94, $c)
    $dst = $ws.Cells.Item(95, $c)
    $src.Copy() | Out-Null
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = $false

# --- Fill in the new row's values ---
$ws.Range("A95").Value = "shexview"
$ws.Range("B95").Value = "x"
$ws.Range("H95").Value = "ShellExView, disable broken context menu entries (e.g. for Win 11)"

# --- Re-apply the AutoFilter so its range grows to include the new last row ---
$lastRow = 129
$ws.AutoFilterMode = $false
$ws.Range("A1:H" + $lastRow).AutoFilter() | Out-Null

# --- Update the hidden _FilterDatabase defined name to match the new range ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$H`$" + $lastRow
    }
}

# --- Restore selection near where the edit happened ---
$ws.Activate() | Out-Null
$ws.Range("H95").Select() | Out-Null
